$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 2.4
$ws.Range("I2").Value = 3.4
$ws.Range("J2").Value = 3.25
$ws.Range("L2").Value = 4
$ws.Range("O2").Value = 1.57
$ws.Range("P2").Value = 2.25
$ws.Range("Q2").Value = 2.88
$ws.Range("R2").Value = 1.4
$ws.Range("S2").Value = 1.62
$ws.Range("T2").Value = 2.2
$ws.Range("U2").Value = 2.25
$ws.Range("V2").Value = 1.57
$ws.Range("W2").Value = 6
$ws.Range("X2").Value = 10
$ws.Range("Z2").Value = 23
$ws.Range("AC2").Value = 5.5
$ws.Range("AE2").Value = 19
$ws.Range("AH2").Value = 15
$ws.Range("AI2").Value = 13
$ws.Range("AK2").Value = 34
$ws.Range("AN2").Value = 4.33
$ws.Range("AT2").Value = 2.2
$ws.Range("AU2").Value = 9.5
$ws.Range("AW2").Value = 5
$ws.Range("AX2").Value = 21
$ws.Range("AY2").Value = 34
$ws.Range("AZ2").Value = 67
$ws.Range("BA2").Value = 126
$ws.Range("Q4").Value = 2.7
$ws.Range("R4").Value = 1.44
$ws.Range("G7").Value = 2.25
$ws.Range("H7").Value = 2.9
$ws.Range("I7").Value = 3.35
$ws.Range("J7").Value = 2.85
$ws.Range("L7").Value = 3.85
$ws.Range("S7").Value = 1.44
$ws.Range("U7").Value = 1.83
$ws.Range("X7").Value = 10
$ws.Range("Y7").Value = 9
$ws.Range("Z7").Value = 23
$ws.Range("AA7").Value = 21
$ws.Range("AC7").Value = 7.2
$ws.Range("AD7").Value = 5.7
$ws.Range("AE7").Value = 15
$ws.Range("AF7").Value = 80
$ws.Range("AI7").Value = 11.75
$ws.Range("AL7").Value = 45
$ws.Range("AN7").Value = 4.05
$ws.Range("AO7").Value = 12
$ws.Range("AQ7").Value = 50
$ws.Range("AT7").Value = 2.37
$ws.Range("AU7").Value = 6.9
$ws.Range("AV7").Value = 65
$ws.Range("AX7").Value = 18.5
$ws.Range("AY7").Value = 25
$ws.Range("AZ7").Value = 100
$ws.Range("BB7").Value = 350
$ws.Range("G8").Value = 2.4
$ws.Range("M8").Value = 1.07
$ws.Range("O8").Value = 1.36
$ws.Range("Q8").Value = 2.2
$ws.Range("R8").Value = 1.65
$ws.Range("S8").Value = 1.44
$ws.Range("T8").Value = 2.63
$ws.Range("AT8").Value = 2.63
$ws.Range("I13").Value = 4.9
$ws.Range("J13").Value = 2.32
$ws.Range("L13").Value = 5.1
$ws.Range("P13").Value = 2.57
$ws.Range("Q13").Value = 2.12
$ws.Range("U13").Value = 1.98
$ws.Range("V13").Value = 1.65
$ws.Range("AC13").Value = 7.5
$ws.Range("AG13").Value = 11.25
$ws.Range("AJ13").Value = 100
$ws.Range("AP13").Value = 19
$ws.Range("AW13").Value = 6.4
$ws.Range("AX13").Value = 29
$ws.Range("AY13").Value = 35
$ws.Range("BA13").Value = 175
$ws.Range("BB13").Value = 500
